$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Update row 4: "дерево" -> "до свидания", drop the sticker id, add answer
$ws1.Range("A4").Value = "до свидания"
$ws1.Range("B4").ClearContents()
$ws1.Range("C4").Value = "и вам не хворать"

# Update row 5: "вадик" -> "до свидания", drop the sticker id, add answer
$ws1.Range("A5").Value = "до свидания"
$ws1.Range("B5").ClearContents()
$ws1.Range("C5").Value = "и вам не хворать"

# Add new "Users" worksheet after Stickers
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Users"

$ws2.Range("A1").Value = "user_id"
$ws2.Range("B1").Value = "имя"
$ws2.Range("C1").Value = "пол"
$ws2.Range("D1").Value = "класс"
$ws2.Range("D1").Select()

# Match page setup used on the rest of the workbook
$ws2.PageSetup.Orientation = 1
$ws2.PageSetup.PaperSize = 9

# Switch back to Stickers sheet and select row 4 (A4:XFD4)
# to match the diff's selection + active tab
$ws1.Activate()
$ws1.Rows.Item(4).Select()
